# Add "% of Revenue from shoes" and "Shoes - Cost of Sales ($M)" columns,
# modeling random per-year inventory holding cost / revenue-mix noise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("C1").Value = '% of Revenue from shoes'
$ws.Range("D1").Value = 'Shoes - Cost of Sales ($M)'

# Random-looking percentages (column C) per year, and the derived
# Cost of Sales for shoes (column D = B*C).
$pct = @{
    2 = 0.68032029565753005
    3 = 0.65735101727793566
    4 = 0.66293650042585406
    5 = 0.65577691485170808
    6 = 0.65154938670109752
    7 = 0.64738203913131964
    8 = 0.65550373134328355
    9 = 0.65292107511336006
}

for ($r = 2; $r -le 9; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $pct[$r]
    $cCell.Style = "Percent"
    $cCell.NumberFormat = "0%"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Formula = "=B$r*C$r"
    $dCell.Style = "Comma"
    $dCell.NumberFormat = '"$"#,##0.00'
}

$ws.Range("E10").Select()
